$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2140350877192982
$ws.Range("C2").Value = 0.5298245614035088
$ws.Range("J2").Value = 0.01052631578947368
$ws.Range("P2").Value = 0.1403508771929824
$ws.Range("S2").Value = 0.1052631578947368
$ws.Range("B3").Value = 0.0189873417721519
$ws.Range("C3").Value = 0.0379746835443038
$ws.Range("J3").Value = 0.01265822784810127
$ws.Range("P3").Value = 0.7531645569620253
$ws.Range("S3").Value = 0.1772151898734177
$ws.Range("J4").Value = 0.08333333333333333
$ws.Range("P4").Value = 0.7777777777777778
$ws.Range("S4").Value = 0.1388888888888889
$ws.Range("B6").Value = 0.095
$ws.Range("D6").Value = 0.005
$ws.Range("F6").Value = 0.06
$ws.Range("J6").Value = 0.305
$ws.Range("Q6").Value = 0.12
$ws.Range("R6").Value = 0.045
$ws.Range("S6").Value = 0.37
$ws.Range("B7").Value = 0.08870967741935484
$ws.Range("D7").Value = 0.02016129032258064
$ws.Range("F7").Value = 0.03225806451612903
$ws.Range("J7").Value = 0.1733870967741936
$ws.Range("O7").Value = 0.01612903225806452
$ws.Range("Q7").Value = 0.1209677419354839
$ws.Range("R7").Value = 0.1209677419354839
$ws.Range("S7").Value = 0.4274193548387097
$ws.Range("B8").Value = 0.07428571428571429
$ws.Range("D8").Value = 0.009523809523809525
$ws.Range("E8").Value = 0.001904761904761905
$ws.Range("F8").Value = 0.0419047619047619
$ws.Range("J8").Value = 0.09333333333333334
$ws.Range("O8").Value = 0.02285714285714286
$ws.Range("Q8").Value = 0.1371428571428571
$ws.Range("R8").Value = 0.1352380952380952
$ws.Range("S8").Value = 0.4838095238095238
$ws.Range("B9").Value = 0.1011904761904762
$ws.Range("D9").Value = 0.01785714285714286
$ws.Range("F9").Value = 0.04761904761904762
$ws.Range("J9").Value = 0.08333333333333333
$ws.Range("O9").Value = 0.005952380952380952
$ws.Range("Q9").Value = 0.1428571428571428
$ws.Range("R9").Value = 0.119047619047619
$ws.Range("S9").Value = 0.4821428571428572
$ws.Range("B10").Value = 0.09147286821705426
$ws.Range("D10").Value = 0.01705426356589147
$ws.Range("E10").Value = 0.0007751937984496124
$ws.Range("F10").Value = 0.0565891472868217
$ws.Range("J10").Value = 0.1054263565891473
$ws.Range("O10").Value = 0.01085271317829457
$ws.Range("Q10").Value = 0.2201550387596899
$ws.Range("R10").Value = 0.09922480620155039
$ws.Range("S10").Value = 0.3984496124031008
$ws.Range("G11").Value = 0.1514360313315927
$ws.Range("J11").Value = 0.08093994778067885
$ws.Range("K11").Value = 0.216710182767624
$ws.Range("L11").Value = 0.5300261096605744
$ws.Range("S11").Value = 0.02088772845953003
$ws.Range("F12").Value = 0.004807692307692308
$ws.Range("G12").Value = 0.7740384615384616
$ws.Range("J12").Value = 0.1730769230769231
$ws.Range("K12").Value = 0.009615384615384616
$ws.Range("L12").Value = 0.02884615384615385
$ws.Range("S12").Value = 0.009615384615384616
$ws.Range("G13").Value = 0.7727272727272727
$ws.Range("J13").Value = 0.1818181818181818
$ws.Range("S13").Value = 0.04545454545454546
$ws.Range("F15").Value = 0.03153153153153153
$ws.Range("H15").Value = 0.1891891891891892
$ws.Range("I15").Value = 0.07657657657657657
$ws.Range("J15").Value = 0.2927927927927928
$ws.Range("K15").Value = 0.07657657657657657
$ws.Range("M15").Value = 0.01351351351351351
$ws.Range("N15").Value = 0.004504504504504504
$ws.Range("O15").Value = 0.1036036036036036
$ws.Range("S15").Value = 0.2117117117117117
$ws.Range("F16").Value = 0.03278688524590164
$ws.Range("H16").Value = 0.185792349726776
$ws.Range("I16").Value = 0.1038251366120219
$ws.Range("J16").Value = 0.4316939890710382
$ws.Range("K16").Value = 0.07103825136612021
$ws.Range("M16").Value = 0.01092896174863388
$ws.Range("O16").Value = 0.03825136612021858
$ws.Range("S16").Value = 0.1256830601092896
$ws.Range("F17").Value = 0.01631701631701632
$ws.Range("H17").Value = 0.1888111888111888
$ws.Range("I17").Value = 0.05594405594405594
$ws.Range("J17").Value = 0.4079254079254079
$ws.Range("K17").Value = 0.1142191142191142
$ws.Range("M17").Value = 0.02331002331002331
$ws.Range("O17").Value = 0.07925407925407925
$ws.Range("S17").Value = 0.1142191142191142
$ws.Range("F18").Value = 0.02713178294573643
$ws.Range("H18").Value = 0.1666666666666667
$ws.Range("I18").Value = 0.06201550387596899
$ws.Range("J18").Value = 0.4534883720930232
$ws.Range("K18").Value = 0.124031007751938
$ws.Range("M18").Value = 0.003875968992248062
$ws.Range("O18").Value = 0.06976744186046512
$ws.Range("S18").Value = 0.09302325581395349
$ws.Range("F19").Value = 0.01564722617354196
$ws.Range("H19").Value = 0.2339971550497866
$ws.Range("I19").Value = 0.06756756756756757
$ws.Range("J19").Value = 0.3477951635846372
$ws.Range("K19").Value = 0.131578947368421
$ws.Range("M19").Value = 0.02062588904694168
$ws.Range("N19").Value = 0.001422475106685633
$ws.Range("O19").Value = 0.0561877667140825
$ws.Range("S19").Value = 0.1251778093883357
